$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 14.1346542613174
$ws.Range("C2").Value = 9.817118980061514
$ws.Range("D2").Value = 9.849515524223582
$ws.Range("F2").Value = 33.90844441242037
$ws.Range("G2").Value = 35.81312204853899
$ws.Range("H2").Value = 15.8530066540732
$ws.Range("I2").Value = 22.14315479450082
$ws.Range("J2").Value = 11.00411746557383
$ws.Range("M2").Value = 17.95394511809635
$ws.Range("B3").Value = 13.57018978290471
$ws.Range("C3").Value = 9.308605990471516
$ws.Range("D3").Value = 9.837365958358838
$ws.Range("F3").Value = 33.93766124023617
$ws.Range("G3").Value = 35.76469206226051
$ws.Range("H3").Value = 15.90644387822675
$ws.Range("I3").Value = 22.26834258116095
$ws.Range("J3").Value = 11.03497496351174
$ws.Range("M3").Value = 17.77782849845438
$ws.Range("B4").Value = 13.21309751312821
$ws.Range("C4").Value = 8.98278789639039
$ws.Range("D4").Value = 9.83098024834908
$ws.Range("F4").Value = 33.96667066621266
$ws.Range("G4").Value = 35.75080064459137
$ws.Range("H4").Value = 15.94315189030491
$ws.Range("I4").Value = 22.35145930475746
$ws.Range("J4").Value = 11.05547349071948
$ws.Range("M4").Value = 17.67187964615443
$ws.Range("B5").Value = 13.0651526369273
$ws.Range("C5").Value = 8.846746600804803
$ws.Range("D5").Value = 9.828649922582311
$ws.Range("F5").Value = 33.98126559098891
$ws.Range("G5").Value = 35.74911495849526
$ws.Range("H5").Value = 15.95908713361269
$ws.Range("I5").Value = 22.38689553471874
$ws.Range("J5").Value = 11.06421704600824
$ws.Range("M5").Value = 17.62929584763052
$ws.Range("B6").Value = 13.04044672080372
$ws.Range("C6").Value = 8.823964570120742
$ws.Range("D6").Value = 9.828279446140874
$ws.Range("F6").Value = 33.98385623255423
$ws.Range("G6").Value = 35.7490747692278
$ws.Range("H6").Value = 15.96179203565898
$ws.Range("I6").Value = 22.3928740610288
$ws.Range("J6").Value = 11.06569247626391
$ws.Range("M6").Value = 17.62226173006941
$ws.Range("B7").Value = 13.21111179975157
$ws.Range("C7").Value = 8.980966210207942
$ws.Range("D7").Value = 9.830947717547762
$ws.Range("F7").Value = 33.96685628440991
$ws.Range("G7").Value = 35.7507618317189
$ws.Range("H7").Value = 15.94336284990934
$ws.Range("I7").Value = 22.3519308799571
$ws.Range("J7").Value = 11.05558982930502
$ws.Range("M7").Value = 17.67130289877173
$ws.Range("B8").Value = 13.9423297233175
$ws.Range("C8").Value = 9.644685770468801
$ws.Range("D8").Value = 9.845104516979875
$ws.Range("F8").Value = 33.91621509198172
$ws.Range("G8").Value = 35.79313098206138
$ws.Range("H8").Value = 15.87062079021742
$ws.Range("I8").Value = 22.18501832988871
$ws.Range("J8").Value = 11.01443503581181
$ws.Range("M8").Value = 17.89279110547861
$ws.Range("B9").Value = 15.2846214010293
$ws.Range("C9").Value = 10.83308611737706
$ws.Range("D9").Value = 9.881306747874541
$ws.Range("F9").Value = 33.9051232572338
$ws.Range("G9").Value = 36.00208180551883
$ws.Range("H9").Value = 15.75905375775571
$ws.Range("I9").Value = 21.90758780399759
$ws.Range("J9").Value = 10.94604424025024
$ws.Range("M9").Value = 18.34263306741808
$ws.Range("B10").Value = 16.20550423568261
$ws.Range("C10").Value = 11.63153617743696
$ws.Range("D10").Value = 9.912934839443523
$ws.Range("F10").Value = 33.95114543164254
$ws.Range("G10").Value = 36.23214369188554
$ws.Range("H10").Value = 15.69624003525115
$ws.Range("I10").Value = 21.73456576966748
$ws.Range("J10").Value = 10.90330435764936
$ws.Range("M10").Value = 18.68005868544091
$ws.Range("B11").Value = 16.6086750777298
$ws.Range("C11").Value = 11.97769883935496
$ws.Range("D11").Value = 9.928388647657741
$ws.Range("F11").Value = 33.98388615644146
$ws.Range("G11").Value = 36.35326631209539
$ws.Range("H11").Value = 15.67186555407122
$ws.Range("I11").Value = 21.66262926390459
$ws.Range("J11").Value = 10.88549102353593
$ws.Range("M11").Value = 18.83452314942635
$ws.Range("B12").Value = 16.7589734285528
$ws.Range("C12").Value = 12.10627484857927
$ws.Range("D12").Value = 9.93439138915863
$ws.Range("F12").Value = 33.99798155345935
$ws.Range("G12").Value = 36.401477841032
$ws.Range("H12").Value = 15.66324271405983
$ws.Range("I12").Value = 21.63636994701476
$ws.Range("J12").Value = 10.8789798999627
$ws.Range("M12").Value = 18.89310847369052
$ws.Range("B13").Value = 16.72671114229334
$ws.Range("C13").Value = 12.07869603025356
$ws.Range("D13").Value = 9.933091930277882
$ws.Range("F13").Value = 33.99487040327732
$ws.Range("G13").Value = 36.390990752165
$ws.Range("H13").Value = 15.66507273827037
$ws.Range("I13").Value = 21.6419815866115
$ws.Range("J13").Value = 10.88037176011217
$ws.Range("M13").Value = 18.88048767133622
$ws.Range("B14").Value = 16.62108835904571
$ws.Range("C14").Value = 11.98832741417196
$ws.Range("D14").Value = 9.92887949131581
$ws.Range("F14").Value = 33.98501178875232
$ws.Range("G14").Value = 36.35718585620295
$ws.Range("H14").Value = 15.6711439592455
$ws.Range("I14").Value = 21.66044917973765
$ws.Range("J14").Value = 10.88495065110855
$ws.Range("M14").Value = 18.83934139886363
$ws.Range("B15").Value = 16.55607923949176
$ws.Range("C15").Value = 11.93264591369047
$ws.Range("D15").Value = 9.926318800628801
$ws.Range("F15").Value = 33.97919407963606
$ws.Range("G15").Value = 36.33678400074253
$ws.Range("H15").Value = 15.6749419347284
$ws.Range("I15").Value = 21.67188915946642
$ws.Range("J15").Value = 10.88778588319617
$ws.Range("M15").Value = 18.81414886424913
$ws.Range("B16").Value = 16.17882952267372
$ws.Range("C16").Value = 11.60856540694052
$ws.Range("D16").Value = 9.911946130571831
$ws.Range("F16").Value = 33.94924331542626
$ws.Range("G16").Value = 36.22455756489332
$ws.Range("H16").Value = 15.69791776309273
$ws.Range("I16").Value = 21.73940391428858
$ws.Range("J16").Value = 10.90450129547423
$ws.Range("M16").Value = 18.66997958790981
$ws.Range("B17").Value = 15.94328350360599
$ws.Range("C17").Value = 11.40534078898014
$ws.Range("D17").Value = 9.903400204151373
$ws.Range("F17").Value = 33.93389348531238
$ws.Range("G17").Value = 36.15991323984422
$ws.Range("H17").Value = 15.71309077137836
$ws.Range("I17").Value = 21.78256195038962
$ws.Range("J17").Value = 10.9151730319026
$ws.Range("M17").Value = 18.58175132061671
$ws.Range("B18").Value = 15.80632730432712
$ws.Range("C18").Value = 11.28684848471406
$ws.Range("D18").Value = 9.89858531459849
$ws.Range("F18").Value = 33.9261759720576
$ws.Range("G18").Value = 36.12428314628676
$ws.Range("H18").Value = 15.72221300006978
$ws.Range("I18").Value = 21.80802226918811
$ws.Range("J18").Value = 10.92146447850872
$ws.Range("M18").Value = 18.5310986019252
$ws.Range("B19").Value = 15.7597063324805
$ws.Range("C19").Value = 11.24645561077774
$ws.Range("D19").Value = 9.896972413461617
$ws.Range("F19").Value = 33.92375379683025
$ws.Range("G19").Value = 36.11248650369208
$ws.Range("H19").Value = 15.72536938772716
$ws.Range("I19").Value = 21.81675186166651
$ws.Range("J19").Value = 10.92362098807695
$ws.Range("M19").Value = 18.51396599349497
$ws.Range("B20").Value = 15.96851142306416
$ws.Range("C20").Value = 11.42714072665158
$ws.Range("D20").Value = 9.904299549391313
$ws.Range("F20").Value = 33.93541248070952
$ws.Range("G20").Value = 36.16663430005913
$ws.Range("H20").Value = 15.71143466114771
$ws.Range("I20").Value = 21.77790172568553
$ws.Range("J20").Value = 10.91402113575952
$ws.Range("M20").Value = 18.59113399329171
$ws.Range("B21").Value = 16.65217754947883
$ws.Range("C21").Value = 12.01493933000963
$ws.Range("D21").Value = 9.930112717341361
$ws.Range("F21").Value = 33.98786145512413
$ws.Range("G21").Value = 36.36705173959726
$ws.Range("H21").Value = 15.66934418828712
$ws.Range("I21").Value = 21.65499810160055
$ws.Range("J21").Value = 10.88359935737831
$ws.Range("M21").Value = 18.8514248989727
$ws.Range("B22").Value = 17.08511696910556
$ws.Range("C22").Value = 12.38446112241091
$ws.Range("D22").Value = 9.947860350796907
$ws.Range("F22").Value = 34.03203072633508
$ws.Range("G22").Value = 36.51169026106578
$ws.Range("H22").Value = 15.6453765576354
$ws.Range("I22").Value = 21.58039787611975
$ws.Range("J22").Value = 10.8650832793675
$ws.Range("M22").Value = 19.02206095210112
$ws.Range("B23").Value = 16.8553509920908
$ws.Range("C23").Value = 12.18859543082521
$ws.Range("D23").Value = 9.938308717386843
$ws.Range("F23").Value = 34.00755242780274
$ws.Range("G23").Value = 36.43325356080074
$ws.Range("H23").Value = 15.65784348265322
$ws.Range("I23").Value = 21.61968708990182
$ws.Range("J23").Value = 10.87484060298684
$ws.Range("M23").Value = 18.93095675817493
$ws.Range("B24").Value = 15.95711066163492
$ws.Range("C24").Value = 11.41729012792304
$ws.Range("D24").Value = 9.903892649104415
$ws.Range("F24").Value = 33.93472229357
$ws.Range("G24").Value = 36.16359092617334
$ws.Range("H24").Value = 15.71218214569582
$ws.Range("I24").Value = 21.78000659466381
$ws.Range("J24").Value = 10.9145414216889
$ws.Range("M24").Value = 18.5868918626265
$ws.Range("B25").Value = 14.93236260428664
$ws.Range("C25").Value = 10.52436391283673
$ws.Range("D25").Value = 9.870620875216634
$ws.Range("F25").Value = 33.89863119153377
$ws.Range("G25").Value = 35.93207680214131
$ws.Range("H25").Value = 15.78588635569415
$ws.Range("I25").Value = 21.97725912591411
$ws.Range("J25").Value = 10.9632274972867
$ws.Range("M25").Value = 18.21953807310407
